$wb = $excel.ActiveWorkbook

# zh-cn sheet (row 2 -> new handback report timestamps for 4f3491b5... file)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-23 20:55:06"
$wsZh.Range("H2").Value = "2016-03-23 20:55:33"

# de-de sheet (row 2 -> new handback report timestamps for 4f3491b5... file)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-23 20:55:11"
$wsDe.Range("H2").Value = "2016-03-23 20:55:40"
